## Update environmental data: Lightstation data / PDO data / Site slope and area
## Applies the Sheet1 restructuring: sort existing slope/aspect rows, add a
## "compass" column with qualitative direction/shallowest-steepest labels,
## append per-site "mean" summary rows (AVERAGE of site + transect scales),
## apply a 0.00 number format to all numeric value cells, and turn the
## existing data range into a filtered/sorted table the same way Sheet2 already is.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

## 1. Re-sort the existing A1:E13 table by scale (C) then by measure (B),
##    same ordering Sheet2's filtered table already uses.
$ws1.Range("A1:E13").Sort($ws1.Range("C2:C13"), 1, $ws1.Range("B2:B13"), [Type]::Missing, 1)

## 2. Turn the (still 13-row) data range into a filtered range, same as
##    Sheet2 already is, before the summary rows below extend the used
##    range - so the filter stays scoped to the original table.
$ws1.Range("A1:E13").AutoFilter()

## 3. New "compass" column (F) - qualitative direction for the "aspect"
##    site rows, and shallowest/steepest markers for the "slope" site rows.
$ws1.Range("F1").Value2 = "compass"
$ws1.Range("F2").Value2 = "SbW"
$ws1.Range("F3").Value2 = "SSW"
$ws1.Range("F4").Value2 = "SSW"
$ws1.Range("F5").Value2 = "shallowest"
$ws1.Range("F7").Value2 = "steepest"

## 4. Append per-site "mean" rows: average of the matching "site" and
##    "transect" rows for each site/measure combination.
$ws1.Range("A14").Value2 = "foggy_cove"
$ws1.Range("B14").Value2 = "aspect"
$ws1.Range("C14").Value2 = "mean"
$ws1.Range("D14").Formula = "=AVERAGE(D2,D8)"

$ws1.Range("A15").Value2 = "north_beach"
$ws1.Range("B15").Value2 = "aspect"
$ws1.Range("C15").Value2 = "mean"
$ws1.Range("D15").Formula = "=AVERAGE(D3,D9)"

$ws1.Range("A16").Value2 = "fifth_beach"
$ws1.Range("B16").Value2 = "aspect"
$ws1.Range("C16").Value2 = "mean"
$ws1.Range("D16").Formula = "=AVERAGE(D4,D10)"

$ws1.Range("A17").Value2 = "foggy_cove"
$ws1.Range("B17").Value2 = "slope"
$ws1.Range("C17").Value2 = "mean"
$ws1.Range("D17").Formula = "=AVERAGE(D5,D11)"

$ws1.Range("A18").Value2 = "north_beach"
$ws1.Range("B18").Value2 = "slope"
$ws1.Range("C18").Value2 = "mean"
$ws1.Range("D18").Formula = "=AVERAGE(D6,D12)"

$ws1.Range("A19").Value2 = "fifth_beach"
$ws1.Range("B19").Value2 = "slope"
$ws1.Range("C19").Value2 = "mean"
$ws1.Range("D19").Formula = "=AVERAGE(D7,D13)"

## 5. Number formatting: mean / standard_deviation columns to 2 decimals.
$ws1.Range("D2:E13").NumberFormat = "0.00"
$ws1.Range("D14:D19").NumberFormat = "0.00"

## 6. Column widths, matching the wider "compass" layout.
$ws1.Columns.Item(1).ColumnWidth = 17.6
$ws1.Columns.Item(2).ColumnWidth = 15.8
$ws1.Columns.Item(4).ColumnWidth = 17.5
$ws1.Columns.Item(5).ColumnWidth = 17.9

## 7. Register the (now filtered) Sheet1 range as its own hidden
##    _FilterDatabase defined name, same as Sheet2 already has - rebuilt in
##    sheet order so Sheet1's entry precedes Sheet2's.
$existingFilterNames = @()
foreach ($n in $wb.Names) {
  if ($n.Name -like "*_FilterDatabase*") { $existingFilterNames += $n }
}
foreach ($n in $existingFilterNames) { $n.Delete() }
$ws1.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$13")
$ws2.Names.Add("_xlnm._FilterDatabase", "=Sheet2!`$A`$1:`$E`$13")

## 8. Select the cell the author left active, matching the saved view state.
$ws1.Range("H13").Select()
